# Auto-generated Excel COM-interop script
# Applies the GitHub Actions crypto-price refresh described in the commit:
#   "Updated cryptos list on Tue Apr  9 23:44:50 UTC 2024 with GitHub Actions"
#
# Updates Price (column D) and Volume(1h) (column E) for the existing rows,
# and swaps the Monero / MXToken rows (48/49 -> rows 50/51) with their refreshed
# price + link + volume data.
#
# Values are written with a leading apostrophe (forces literal text, the way
# typing into a cell formatted as this sheet already stores Price/Volume: as
# text, e.g. '68.948.59' - a thousand-separated string, not a real number) and
# then the cell style is reset to Normal so the text-coercion doesn't leave a
# stray numeric/quote-prefix style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '69.024.22' }
    @{ Cell = 'E2'; Value = '  -3.80%  ' }
    @{ Cell = 'D3'; Value = '3.495.02' }
    @{ Cell = 'E3'; Value = '  -5.48%  ' }
    @{ Cell = 'E4'; Value = '  -0.10%  ' }
    @{ Cell = 'D5'; Value = '579.85' }
    @{ Cell = 'E5'; Value = '  -1.36%  ' }
    @{ Cell = 'D6'; Value = '172.48' }
    @{ Cell = 'E6'; Value = '  -4.56%  ' }
    @{ Cell = 'E7'; Value = '  -0.61%  ' }
    @{ Cell = 'D8'; Value = '3.489.03' }
    @{ Cell = 'E8'; Value = '  -5.36%  ' }
    @{ Cell = 'E9'; Value = '  -0.05%  ' }
    @{ Cell = 'D10'; Value = '0.188' }
    @{ Cell = 'E10'; Value = '  -7.35%  ' }
    @{ Cell = 'D11'; Value = '6.67' }
    @{ Cell = 'E11'; Value = '  +4.33%  ' }
    @{ Cell = 'E12'; Value = '  -3.68%  ' }
    @{ Cell = 'D13'; Value = '46.44' }
    @{ Cell = 'E13'; Value = '  -6.87%  ' }
    @{ Cell = 'E14'; Value = '  -4.83%  ' }
    @{ Cell = 'D15'; Value = '671.80' }
    @{ Cell = 'E15'; Value = '  -1.46%  ' }
    @{ Cell = 'D16'; Value = '4.052.47' }
    @{ Cell = 'E16'; Value = '  -5.64%  ' }
    @{ Cell = 'D17'; Value = '8.65' }
    @{ Cell = 'E17'; Value = '  -4.35%  ' }
    @{ Cell = 'D18'; Value = '68.937.92' }
    @{ Cell = 'E18'; Value = '  -4.03%  ' }
    @{ Cell = 'D19'; Value = '3.491.33' }
    @{ Cell = 'E19'; Value = '  -5.68%  ' }
    @{ Cell = 'E20'; Value = '  -1.46%  ' }
    @{ Cell = 'D21'; Value = '17.32' }
    @{ Cell = 'E21'; Value = '  -4.35%  ' }
    @{ Cell = 'D22'; Value = '11.12' }
    @{ Cell = 'E22'; Value = '  -4.49%  ' }
    @{ Cell = 'D23'; Value = '0.896' }
    @{ Cell = 'E23'; Value = '  -5.36%  ' }
    @{ Cell = 'D24'; Value = '16.02' }
    @{ Cell = 'E24'; Value = '  -10.08%  ' }
    @{ Cell = 'D25'; Value = '97.34' }
    @{ Cell = 'E25'; Value = '  -5.87%  ' }
    @{ Cell = 'E26'; Value = '  -5.03%  ' }
    @{ Cell = 'E27'; Value = '  +0.02%  ' }
    @{ Cell = 'D28'; Value = '2.64' }
    @{ Cell = 'E28'; Value = '  -7.18%  ' }
    @{ Cell = 'D29'; Value = '9.36' }
    @{ Cell = 'E29'; Value = '  -9.51%  ' }
    @{ Cell = 'D30'; Value = '32.73' }
    @{ Cell = 'E30'; Value = '  -8.05%  ' }
    @{ Cell = 'D31'; Value = '8.63' }
    @{ Cell = 'E31'; Value = '  -7.41%  ' }
    @{ Cell = 'D32'; Value = '3.17' }
    @{ Cell = 'E32'; Value = '  -8.97%  ' }
    @{ Cell = 'D33'; Value = '1.36' }
    @{ Cell = 'E33'; Value = '  -5.96%  ' }
    @{ Cell = 'D34'; Value = '7.21' }
    @{ Cell = 'E34'; Value = '  -2.01%  ' }
    @{ Cell = 'D35'; Value = '591.74' }
    @{ Cell = 'E35'; Value = '  +4.53%  ' }
    @{ Cell = 'D36'; Value = '10.81' }
    @{ Cell = 'E36'; Value = '  -4.31%  ' }
    @{ Cell = 'D37'; Value = '3.56' }
    @{ Cell = 'E37'; Value = '  -15.13%  ' }
    @{ Cell = 'E38'; Value = '  -5.61%  ' }
    @{ Cell = 'D39'; Value = '56.88' }
    @{ Cell = 'E39'; Value = '  -4.41%  ' }
    @{ Cell = 'D40'; Value = '0.999' }
    @{ Cell = 'E40'; Value = '  +0.00%  ' }
    @{ Cell = 'D41'; Value = '0.0435' }
    @{ Cell = 'E41'; Value = '  -6.48%  ' }
    @{ Cell = 'D42'; Value = '0.333' }
    @{ Cell = 'E42'; Value = '  -5.59%  ' }
    @{ Cell = 'D43'; Value = '3.397.94' }
    @{ Cell = 'E43'; Value = '  -9.64%  ' }
    @{ Cell = 'E44'; Value = '  -6.62%  ' }
    @{ Cell = 'D45'; Value = '33.05' }
    @{ Cell = 'E45'; Value = '  -7.41%  ' }
    @{ Cell = 'D46'; Value = '0.0₃0702' }
    @{ Cell = 'E46'; Value = '  -9.77%  ' }
    @{ Cell = 'E47'; Value = '  -1.04%  ' }
    @{ Cell = 'D48'; Value = '2.58' }
    @{ Cell = 'E48'; Value = '  -8.30%  ' }
    @{ Cell = 'E49'; Value = '  -1.27%  ' }
    @{ Cell = 'B50'; Value = 'Monero' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' }
    @{ Cell = 'D50'; Value = '132.78' }
    @{ Cell = 'E50'; Value = '  -2.28%  ' }
    @{ Cell = 'B51'; Value = 'MXToken' }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = 'D51'; Value = '5.68' }
    @{ Cell = 'E51'; Value = '  +15.31%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $range.Value = "'" + $u.Value
    $range.Style = 'Normal'
}

